# Generate Report for Handoff
# The "aec04b78-29d6-4ae3-98cb-2cbae4683c17.md" row's latest handoff/handback
# timestamps were refreshed on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the
# aec04b78-... row (row 6).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-09-02 02:48:34"

# zh-cn sheet: "Latest Handoff Datetime" column (H) for the aec04b78-... row
# (row 6).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-09-02 02:48:29"

# de-de sheet: "Latest Handoff Datetime" column (H) for the aec04b78-... row
# (row 6).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-09-02 02:48:34"
